$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20 corresponds to S.No 12 ("addProduct() method"): lower the score from 10 to 9
# and add a grading comment explaining the deduction.
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "(-1)For not getting products of a customer into linked list."

# Move the active selection to F20 (matches the saved cursor position in the file)
$ws.Range("F20").Select()
